$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 18; $r++) {
    $ws.Range("AO$r").Value = 162709.28305112128
}

$excel.CalculateFullRebuild()
